# Apply the "Added list view for exam schedule preview" changes:
# Electives B6/B7 for Section_A and Section_B are reassigned to different
# rooms, the Classroom_Utilization stats for the affected rooms are
# recomputed, and the Classroom_Allocation sheet is updated to reflect the
# new room numbers / types / capacities.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Section_A
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Section_A")
$ws.Range("B2").Value = "ELECTIVE_B6 [C003]"
$ws.Range("C2").Value = "ELECTIVE_B7 [C404]"
$ws.Range("D5").Value = "ELECTIVE_B6 [C003]"
$ws.Range("E5").Value = "ELECTIVE_B7 [C404]"
$ws.Range("C6").Value = "ELECTIVE_B6 (Tutorial) [C204]"
$ws.Range("D6").Value = "ELECTIVE_B7 (Tutorial) [C003]"

# ---------------------------------------------------------------------
# Section_B
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Section_B")
$ws.Range("B2").Value = "ELECTIVE_B6 [C205]"
$ws.Range("C2").Value = "ELECTIVE_B7 [C101]"
$ws.Range("D5").Value = "ELECTIVE_B6 [C205]"
$ws.Range("E5").Value = "ELECTIVE_B7 [C101]"
$ws.Range("C6").Value = "ELECTIVE_B6 (Tutorial) [C305]"
$ws.Range("D6").Value = "ELECTIVE_B7 (Tutorial) [C201]"

# ---------------------------------------------------------------------
# Classroom_Utilization - update weekly hours / daily avg / utilization
# rate for every room whose weekly schedule changed above.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Classroom_Utilization")

# C003 (row 4)
$ws.Range("D4").Value = 4
$ws.Range("E4").Value = 0.8
$ws.Range("G4").Value = 10

# C102 (row 7)
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("G7").Value = 0

# C201 (row 13)
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 0.2
$ws.Range("G13").Value = 2.5

# C204 (row 16)
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 0.2
$ws.Range("G16").Value = 2.5

# C305 (row 25)
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 0.2
$ws.Range("G25").Value = 2.5

# C401 (row 29)
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("G29").Value = 0

# C404 (row 32)
$ws.Range("D32").Value = 3
$ws.Range("E32").Value = 0.6
$ws.Range("G32").Value = 7.5

# C405 (row 33)
$ws.Range("D33").Value = 0
$ws.Range("E33").Value = 0
$ws.Range("G33").Value = 0

# ---------------------------------------------------------------------
# Classroom_Allocation - mirror the new room numbers/types/capacities.
# Capacity values are stored as text in this sheet, so they are entered
# with a leading apostrophe to keep them as text instead of numbers.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Classroom_Allocation")

$ws.Range("G2").Value = "C003"
$ws.Range("H2").Value = "large classroom"
$ws.Range("I2").Value = "'135"

$ws.Range("G3").Value = "C404"
$ws.Range("I3").Value = "'78"

$ws.Range("G4").Value = "C204"
$ws.Range("I4").Value = "'96"

$ws.Range("G5").Value = "C003"
$ws.Range("H5").Value = "large classroom"
$ws.Range("I5").Value = "'135"

$ws.Range("G6").Value = "C003"
$ws.Range("H6").Value = "large classroom"
$ws.Range("I6").Value = "'135"

$ws.Range("G7").Value = "C404"
$ws.Range("I7").Value = "'78"

$ws.Range("G8").Value = "C205"

$ws.Range("G9").Value = "C101"

$ws.Range("G10").Value = "C305"

$ws.Range("G11").Value = "C205"

$ws.Range("G12").Value = "C201"
$ws.Range("I12").Value = "'96"

$ws.Range("G13").Value = "C101"
